# Scheduled-runner update: refresh cached Universalis market-price figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job-class
# profit sheets. Values only; no structural/formula changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2754.0244
$ws.Range("J17").Value = 2820.125
$ws.Range("L17").Value = 8460.375
$ws.Range("N17").Value = -8796.375

$ws.Range("H112").Value = 5874.5
$ws.Range("I112").Value = 35250
$ws.Range("J112").Value = 1678
$ws.Range("K112").Value = 105750
$ws.Range("L112").Value = 5034
$ws.Range("M112").Value = -104642
$ws.Range("N112").Value = -7250

$ws.Range("H113").Value = 3025.625
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 3041
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 3041
$ws.Range("M113").Value = 254
$ws.Range("N113").Value = -9549

$ws.Range("H137").Value = 531448.9
$ws.Range("I137").Value = 2601.8845
$ws.Range("J137").Value = 903071.0600000001
$ws.Range("K137").Value = 7805.6535
$ws.Range("L137").Value = 2709213.18
$ws.Range("M137").Value = -5255.6535
$ws.Range("N137").Value = -2714313.18

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1694.1177
$ws.Range("I110").Value = 1714.2858
$ws.Range("J110").Value = 1600
$ws.Range("K110").Value = 1714.2858
$ws.Range("L110").Value = 1600
$ws.Range("M110").Value = 330.7141999999999
$ws.Range("N110").Value = -5690

$ws.Range("H122").Value = 1747.1852
$ws.Range("I122").Value = 1734.4
$ws.Range("J122").Value = 1907
$ws.Range("K122").Value = 5203.200000000001
$ws.Range("L122").Value = 5721
$ws.Range("M122").Value = -2753.200000000001
$ws.Range("N122").Value = -10621

$ws.Range("H132").Value = 3841
$ws.Range("I132").Value = 3677.2632
$ws.Range("J132").Value = 4285.4287
$ws.Range("K132").Value = 11031.7896
$ws.Range("L132").Value = 12856.2861
$ws.Range("M132").Value = -8501.7896
$ws.Range("N132").Value = -17916.2861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 33954.13
$ws.Range("I134").Value = 1855.6923
$ws.Range("J134").Value = 200866
$ws.Range("K134").Value = 5567.0769
$ws.Range("L134").Value = 602598
$ws.Range("M134").Value = -3032.0769
$ws.Range("N134").Value = -607668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 40000
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H31").Value = 581248.8
$ws.Range("I31").Value = 14047.8
$ws.Range("K31").Value = 14047.8
$ws.Range("M31").Value = -13752.8

$ws.Range("H34").Value = 581248.8
$ws.Range("I34").Value = 14047.8
$ws.Range("K34").Value = 14047.8
$ws.Range("M34").Value = -13845.8

$ws.Range("H41").Value = 20065
$ws.Range("J41").Value = 20065
$ws.Range("L41").Value = 20065
$ws.Range("N41").Value = -20921

$ws.Range("H132").Value = 3869.625
$ws.Range("I132").Value = 3586.4348
$ws.Range("J132").Value = 4593.3335
$ws.Range("K132").Value = 10759.3044
$ws.Range("L132").Value = 13780.0005
$ws.Range("M132").Value = -8229.304400000001
$ws.Range("N132").Value = -18840.0005

$ws.Range("H133").Value = 22661.375
$ws.Range("J133").Value = 43665
$ws.Range("L133").Value = 43665
$ws.Range("N133").Value = -48725

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 774498.7
$ws.Range("I4").Value = 1084178.2
$ws.Range("K4").Value = 3252534.6
$ws.Range("M4").Value = -3252422.6

$ws.Range("H68").Value = 123394.72
$ws.Range("I68").Value = 245718.62
$ws.Range("J68").Value = 3673.4468
$ws.Range("K68").Value = 737155.86
$ws.Range("L68").Value = 11020.3404
$ws.Range("M68").Value = -736344.86
$ws.Range("N68").Value = -12642.3404

$ws.Range("H71").Value = 123394.72
$ws.Range("I71").Value = 245718.62
$ws.Range("J71").Value = 3673.4468
$ws.Range("K71").Value = 2211467.58
$ws.Range("L71").Value = 33061.0212
$ws.Range("M71").Value = -2207411.58
$ws.Range("N71").Value = -41173.0212

$ws.Range("H108").Value = 3198
$ws.Range("I108").Value = 971
$ws.Range("J108").Value = 5425
$ws.Range("K108").Value = 2913
$ws.Range("L108").Value = 16275
$ws.Range("M108").Value = -33
$ws.Range("N108").Value = -22035

$ws.Range("H124").Value = 2411.375
$ws.Range("I124").Value = 761.6667
$ws.Range("J124").Value = 3401.2
$ws.Range("K124").Value = 2285.0001
$ws.Range("L124").Value = 10203.6
$ws.Range("M124").Value = 2624.9999
$ws.Range("N124").Value = -20023.6

$ws.Range("H130").Value = 2661.75
$ws.Range("I130").Value = 2623.3333
$ws.Range("J130").Value = 2777
$ws.Range("K130").Value = 7869.999899999999
$ws.Range("L130").Value = 8331
$ws.Range("M130").Value = -2849.999899999999
$ws.Range("N130").Value = -18371

$ws.Range("H131").Value = 37865.383
$ws.Range("I131").Value = 1434.9286
$ws.Range("J131").Value = 80367.586
$ws.Range("K131").Value = 4304.7858
$ws.Range("L131").Value = 241102.758
$ws.Range("M131").Value = 735.2142000000003
$ws.Range("N131").Value = -251182.758

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 14100
$ws.Range("I5").Value = 9000
$ws.Range("J5").Value = 17500
$ws.Range("K5").Value = 9000
$ws.Range("L5").Value = 17500
$ws.Range("M5").Value = -8888
$ws.Range("N5").Value = -17724

$ws.Range("H22").Value = 20000
$ws.Range("I22").Value = 20000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 20000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -19471
$ws.Range("N22").ClearContents()

$ws.Range("H25").Value = 3000
$ws.Range("I25").Value = 3000
$ws.Range("K25").Value = 3000
$ws.Range("M25").Value = -2471

$ws.Range("H113").Value = 2255.7273
$ws.Range("I113").Value = 2200
$ws.Range("J113").Value = 2506.5
$ws.Range("K113").Value = 2200
$ws.Range("L113").Value = 2506.5
$ws.Range("M113").Value = -30
$ws.Range("N113").Value = -6846.5

$ws.Range("H132").Value = 8489.227999999999
$ws.Range("I132").Value = 3504.25
$ws.Range("J132").Value = 14471.2
$ws.Range("K132").Value = 10512.75
$ws.Range("L132").Value = 43413.60000000001
$ws.Range("M132").Value = -7982.75
$ws.Range("N132").Value = -48473.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 47978.855
$ws.Range("J133").Value = 47978.855
$ws.Range("L133").Value = 47978.855
$ws.Range("N133").Value = -53038.855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 52400
$ws.Range("I2").Value = 250000
$ws.Range("K2").Value = 250000
$ws.Range("M2").Value = -249888

$ws.Range("H74").Value = 14894.9
$ws.Range("J74").Value = 15438.777
$ws.Range("L74").Value = 15438.777
$ws.Range("N74").Value = -17310.777

$ws.Range("H77").Value = 14894.9
$ws.Range("J77").Value = 15438.777
$ws.Range("L77").Value = 46316.331
$ws.Range("N77").Value = -55676.331

$ws.Range("H136").Value = 6500.4194
$ws.Range("I136").Value = 6857.207
$ws.Range("K136").Value = 20571.621
$ws.Range("M136").Value = -18021.621
